$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 1.95
$ws.Range("I2").Value = 4.5
$ws.Range("J2").Value = 2.75
$ws.Range("L2").Value = 4.75
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("V2").Value = 1.63
$ws.Range("W2").Value = 6
$ws.Range("Z2").Value = 17
$ws.Range("AC2").Value = 7
$ws.Range("AE2").Value = 17
$ws.Range("AJ2").Value = 15
$ws.Range("AU2").Value = 9
$ws.Range("AZ2").Value = 81
$ws.Range("BA2").Value = 126

$wb.Save()
